$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new wishlist entry as the next row after the existing data.
$ws.Range("A17").Value = "La mujer de mi vida"
$ws.Range("B17").Value = "Carla Guelfenbein"

# Column C ("Editorial") is left blank for this entry, matching the other
# rows (8-16) that already have an empty-but-present Editorial cell. A
# plain empty-string assignment removes the cell entirely, so seed it
# with a leading apostrophe (forces an empty text cell) and then reset
# the style back to Normal to drop the quote-prefix formatting that the
# apostrophe trick leaves behind.
$ws.Range("C17").Value = "'"
$ws.Range("C17").Style = "Normal"
